# [ADDITIONAL SCRAPING] added scraping code for extra browling attributes and excel sheets
#
# 1) Clean up "ODI Batting Extra": drop the empty placeholder cells (C/D/E, and
#    occasionally B) that were written as empty inline strings but never carried data.
# 2) Add a new "ODI Bowling Extra" sheet (directly after "ODI Batting Extra") with the
#    scraped MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) ODI Batting Extra - remove empty placeholder cells
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ODI Batting Extra")

$ws4.Range("C3:E3").ClearContents()
$ws4.Range("C5:E5").ClearContents()
$ws4.Range("B6:E6").ClearContents()
$ws4.Range("B7:E7").ClearContents()
$ws4.Range("C8:E8").ClearContents()
$ws4.Range("E9").ClearContents()
$ws4.Range("B10:E10").ClearContents()
$ws4.Range("B16:E16").ClearContents()
$ws4.Range("B17:E17").ClearContents()
$ws4.Range("B18:E18").ClearContents()
$ws4.Range("C19:E19").ClearContents()
$ws4.Range("C20:E20").ClearContents()
$ws4.Range("C21:E21").ClearContents()

# ---------------------------------------------------------------------------
# 2) Add "ODI Bowling Extra" right after "ODI Batting Extra"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "ODI Bowling Extra"

# Reuse the bold/centered/bordered header style from the sibling sheet.
$ws4.Range("A1:C1").Copy()
$ws5.Range("A1:C1").PasteSpecial(-4122)

$ws5.Range("A1").Value = "MATCH_CODE"
$ws5.Range("B1").Value = "MAIDEN_OVERS"
$ws5.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Force the data columns to Text so numeric-looking codes/percentages
# ("0", "10.00%", match codes, ...) are stored verbatim, not re-parsed as numbers.
$ws5.Range("A2:C21").NumberFormat = "@"

$ws5.Range("A2").Value = "3793"
$ws5.Range("B2").Value = "0"
$ws5.Range("C2").Value = "10.00%"

$ws5.Range("A3").Value = "3826"
$ws5.Range("B3").Value = "0"
$ws5.Range("C3").Value = "20.00%"

$ws5.Range("A4").Value = "3827"
$ws5.Range("B4").Value = "0"
$ws5.Range("C4").Value = "20.00%"

$ws5.Range("A5").Value = "3828"
$ws5.Range("B5").Value = "1"
$ws5.Range("C5").Value = "10.00%"

$ws5.Range("A6").Value = "3865"
$ws5.Range("B6").Value = "0"
$ws5.Range("C6").Value = "10.00%"

$ws5.Range("A7").Value = "3868"

$ws5.Range("A8").Value = "3872"
$ws5.Range("B8").Value = "0"
$ws5.Range("C8").Value = "10.00%"

$ws5.Range("A9").Value = "3883"
$ws5.Range("B9").Value = "0"
$ws5.Range("C9").Value = "30.00%"

$ws5.Range("A10").Value = "3884"

$ws5.Range("A11").Value = "3886"
$ws5.Range("B11").Value = "0"

$ws5.Range("A12").Value = "3888"
$ws5.Range("B12").Value = "0"
$ws5.Range("C12").Value = "10.00%"

$ws5.Range("A13").Value = "4026"
$ws5.Range("B13").Value = "1"
$ws5.Range("C13").Value = "10.00%"

$ws5.Range("A14").Value = "4032"
$ws5.Range("B14").Value = "0"
$ws5.Range("C14").Value = "20.00%"

$ws5.Range("A15").Value = "4036"
$ws5.Range("B15").Value = "0"
$ws5.Range("C15").Value = "30.00%"

$ws5.Range("A16").Value = "4039"

$ws5.Range("A17").Value = "4085"

$ws5.Range("A18").Value = "4088"

$ws5.Range("A19").Value = "4089"
$ws5.Range("B19").Value = "0"
$ws5.Range("C19").Value = "20.00%"

$ws5.Range("A20").Value = "4669"
$ws5.Range("B20").Value = "0"
$ws5.Range("C20").Value = "10.00%"

$ws5.Range("A21").Value = "4676"
$ws5.Range("B21").Value = "0"
$ws5.Range("C21").Value = "30.00%"
